# Add a new "Croatia" market test-data sheet, cloned from the existing
# "Slovakia" sheet (same layout/values), placed after "Spain" as the new
# last tab, and make it the active sheet with A12 selected.

$wb = $excel.ActiveWorkbook

$template = $wb.Worksheets.Item("Slovakia")
$lastSheet = $wb.Worksheets.Item("Spain")

# Copy the template sheet so it lands immediately after "Spain" (the
# current last sheet), becoming the new last tab.
$template.Copy($null, $lastSheet)

$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "Croatia"

# Make it the active sheet/tab with A12 selected, matching the target
# workbook view state.
$newSheet.Activate() | Out-Null
$newSheet.Range("A12").Select() | Out-Null
